# Update matchDay.xlsx: shift match dates from round 8 to round 9
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @{
    2  = 45780
    3  = 45781
    4  = 45780
    5  = 45780
    6  = 45778
    7  = 45779
    8  = 45781
    9  = 45780
    10 = 45781
    11 = 45781
    12 = 45780
    13 = 45778
    14 = 45780
    15 = 45781
    16 = 45779
    17 = 45781
    18 = 45780
    19 = 45780
}

foreach ($row in $newDates.Keys) {
    $ws.Range("A$row").Value = $newDates[$row]
}
